$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "VLSM"
$ws.Range("A1").Value = "Hello"
